# Update model assets with corrected translations
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VotingEnsemble")

# Update the "Run ID" value (row 2 in sheet, shared string table entry)
$ws.Range("B2").Value = "e91c8e84-e7d7-471d-83bd-f6628178a777_12"

# Update the selection on the active sheet
$ws.Range("A3:A23").Select()

# Update the metric values in column B that changed
$ws.Range("B3").Value = 0.82857000000000003
$ws.Range("B5").Value = 0.85714000000000001
$ws.Range("B8").Value = 0.83455999999999997
$ws.Range("B10").Value = 0.70833000000000002
$ws.Range("B11").Value = 0.68874000000000002
$ws.Range("B12").Value = 0.82857000000000003
$ws.Range("B13").Value = 0.81299999999999994
$ws.Range("B14").Value = 0.50134000000000001
$ws.Range("B15").Value = 0.42620999999999998
$ws.Range("B17").Value = 0.70404999999999995
$ws.Range("B18").Value = 0.82857000000000003
$ws.Range("B19").Value = 0.83408000000000004
$ws.Range("B20").Value = 0.70833000000000002
$ws.Range("B21").Value = 0.82857000000000003
$ws.Range("B22").Value = 0.82857000000000003
$ws.Range("B23").Value = 0.86075999999999997
